{"js": "const newTexts = [\n  \"Topic: Rust: A Memory-Safe Compiled Programming Language\",\n  \"\",\n  \"Rust is a modern, compiled programming language designed for performance and safety. It is known for its ability to provide low-level control over hardware resources while simultaneously ensuring memory safety without relying on a garbage collector. This makes Rust an excellent choice for system programming, embedded systems, game development, and other areas where performance and reliability are critical. It emphasizes zero-cost abstractions, meaning that the language features do not introduce runtime overhead. Rust's strong emphasis on compile-time checks prevents many common programming errors, leading to more robust and maintainable code.\",\n  \"\",\n  \"*   **Memory Safety:** Rust guarantees memory safety at compile time, preventing dangling pointers and memory leaks.\",\n  \"*   **Performance:** Rust compiles to machine code, resulting in high-performance applications comparable to C and C++.\",\n  \"*   **Concurrency:** Rust offers powerful concurrency features to write concurrent and parallel programs safely.\",\n  \"*   **Ownership and Borrowing:** Rust's unique ownership and borrowing system controls memory access and prevents data races.\",\n  \"*   **Growing Community:** Rust has a vibrant and supportive community, providing a wealth of resources and libraries.\",\n  \"\",\n  \"Topic: High-Performance Applications of Rust\",\n  \"\",\n  \"Rust's design makes it particularly well-suited for high-performance applications where efficiency and control over system resources are paramount. Game engines, such as those used to create 3D games, benefit from Rust's ability to provide fine-grained memory management and optimized code execution. Databases, which require efficient data storage and retrieval, can leverage Rust's performance characteristics to handle large datasets with minimal overhead. Operating systems can be built with Rust's safety features to create robust and secure kernels. Furthermore, Rust is a popular choice for WebAssembly development, enabling high-performance code to run in web browsers.\",\n  \"\",\n  \"*   **Game Engines:** Rust allows for optimization of resource management, leading to smoother gameplay.\",\n  \"*   **Databases:** The language's performance features contribute to efficient data handling and querying.\",\n  \"*   **Operating Systems:** Rust promotes safe and reliable operating system kernels with its memory safety guarantees.\",\n  \"*   **WebAssembly:** Rust enables the development of high-performance web applications and libraries.\",\n  \"*   **Embedded Systems:** It offers low-level control and strong guarantees, and is great for IoT devices.\",\n  \"\",\n  \"Topic: History of Rust\",\n  \"\",\n  \"Rust's journey began as a personal project of Graydon Hoare in 2007. He named the language \\\"Rust\\\" after a type of fungus, reflecting the project's initial state. Mozilla took notice of the language's potential and began sponsoring its development in 2009. This sponsorship significantly accelerated Rust's progress, attracting a dedicated community and expanding its features. Since 2016, Rust has consistently been ranked as one of the most loved programming languages by developers, which is a testament to its impact on the industry. The Rust community is very active, with frequent updates and improvements, making the language constantly evolve.\",\n  \"\",\n  \"*   **Origin:** Started as a personal project by Graydon Hoare in 2007.\",\n  \"*   **Mozilla Sponsorship:** Gained Mozilla's support in 2009, leading to increased development.\",\n  \"*   **Popularity:** Ranked as a top programming language since 2016.\",\n  \"*   **Community:** A large and active community contributes to the language's development and ecosystem.\",\n  \"*   **Evolving:** Rust continually evolves with new features and improvements.\",\n  \"\",\n  \"Topic: Rust's Memory Management Approach\",\n  \"\",\n  \"Rust's approach to memory management differs significantly from that of garbage-collected languages and languages like C and C++. It avoids the overhead of garbage collection while providing the safety of managed memory. Rust achieves this through the concepts of ownership, borrowing, and lifetimes. Ownership dictates which part of the code is responsible for a given piece of memory. Borrowing allows multiple parts of the code to access the same memory, but under strict rules that prevent data races and ensure memory safety. Lifetimes help to track the duration of references, ensuring that they are valid for as long as they are being used.\",\n  \"\",\n  \"*   **Ownership:** Each value has a single owner, ensuring exclusive access.\",\n  \"*   **Borrowing:** Allows multiple parts of the code to access the same memory, but under strict rules to prevent data races.\",\n  \"*   **Lifetimes:** Track the validity of references to ensure they are not used after the memory they point to has been freed.\",\n  \"*   **Compile-Time Safety:** Memory management checks are performed at compile time.\",\n  \"*   **No Garbage Collector:** Ensures low overhead and predictability.\",\n  \"\",\n  \"Topic: Immutability and Memory Allocation in Rust\",\n  \"\",\n  \"By default, variables in Rust are immutable. This means that once a value is assigned to a variable, it cannot be changed. This immutability is a cornerstone of Rust's memory safety, as it eliminates data races that can occur when multiple parts of a program try to modify the same memory location simultaneously. Immutable values, and objects with a known size at compile time, are typically stored on the stack, which is a fast and efficient memory region. Mutable values, as well as objects whose size is not known until runtime, are stored on the heap, a region of memory that is managed dynamically.\",\n  \"\",\n  \"*   **Immutability by Default:** Enhances safety and predictability.\",\n  \"*   **Stack Memory:** Used for immutable values, efficient due to minimal overhead.\",\n  \"*   **Heap Memory:** Used for mutable values and objects with unknown sizes at compile time.\",\n  \"*   **Compile-Time Allocation:** Memory allocation is often determined at compile time, optimizing performance.\",\n  \"*   **Automatic Dropping:** Memory is automatically freed when a variable goes out of scope, preventing memory leaks.\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replace the text of the first 49 paragraphs (indices 0-48) in place,\n// preserving the existing formatting/runs of each paragraph.\nfor (let i = 0; i < newTexts.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// The old document had a trailing \"Topic: Using NPM Modules\" section\n// (10 paragraphs: index 49 blank, 50 topic, 51 blank, 52 body, 53 blank,\n// 54-58 bullets, 59 trailing blank) that the new document no longer has.\n// Paragraph 49 (the blank separator before that section) should remain as\n// the document's final (empty) paragraph, so repeatedly delete paragraph\n// index 49 -- each deletion shifts the following paragraphs up by one,\n// and Word will not delete the very last paragraph of the body, so this\n// converges on removing exactly the 10 obsolete paragraphs.\nfor (let i = 0; i < 10; i++) {\n  const paras = body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  paras.items[49].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$newTexts = @(\n    'Topic: Rust: A Memory-Safe Compiled Programming Language',\n    '',\n    'Rust is a modern, compiled programming language designed for performance and safety. It is known for its ability to provide low-level control over hardware resources while simultaneously ensuring memory safety without relying on a garbage collector. This makes Rust an excellent choice for system programming, embedded systems, game development, and other areas where performance and reliability are critical. It emphasizes zero-cost abstractions, meaning that the language features do not introduce runtime overhead. Rust''s strong emphasis on compile-time checks prevents many common programming errors, leading to more robust and maintainable code.',\n    '',\n    '*   **Memory Safety:** Rust guarantees memory safety at compile time, preventing dangling pointers and memory leaks.',\n    '*   **Performance:** Rust compiles to machine code, resulting in high-performance applications comparable to C and C++.',\n    '*   **Concurrency:** Rust offers powerful concurrency features to write concurrent and parallel programs safely.',\n    '*   **Ownership and Borrowing:** Rust''s unique ownership and borrowing system controls memory access and prevents data races.',\n    '*   **Growing Community:** Rust has a vibrant and supportive community, providing a wealth of resources and libraries.',\n    '',\n    'Topic: High-Performance Applications of Rust',\n    '',\n    'Rust''s design makes it particularly well-suited for high-performance applications where efficiency and control over system resources are paramount. Game engines, such as those used to create 3D games, benefit from Rust''s ability to provide fine-grained memory management and optimized code execution. Databases, which require efficient data storage and retrieval, can leverage Rust''s performance characteristics to handle large datasets with minimal overhead. Operating systems can be built with Rust''s safety features to create robust and secure kernels. Furthermore, Rust is a popular choice for WebAssembly development, enabling high-performance code to run in web browsers.',\n    '',\n    '*   **Game Engines:** Rust allows for optimization of resource management, leading to smoother gameplay.',\n    '*   **Databases:** The language''s performance features contribute to efficient data handling and querying.',\n    '*   **Operating Systems:** Rust promotes safe and reliable operating system kernels with its memory safety guarantees.',\n    '*   **WebAssembly:** Rust enables the development of high-performance web applications and libraries.',\n    '*   **Embedded Systems:** It offers low-level control and strong guarantees, and is great for IoT devices.',\n    '',\n    'Topic: History of Rust',\n    '',\n    'Rust''s journey began as a personal project of Graydon Hoare in 2007. He named the language \"Rust\" after a type of fungus, reflecting the project''s initial state. Mozilla took notice of the language''s potential and began sponsoring its development in 2009. This sponsorship significantly accelerated Rust''s progress, attracting a dedicated community and expanding its features. Since 2016, Rust has consistently been ranked as one of the most loved programming languages by developers, which is a testament to its impact on the industry. The Rust community is very active, with frequent updates and improvements, making the language constantly evolve.',\n    '',\n    '*   **Origin:** Started as a personal project by Graydon Hoare in 2007.',\n    '*   **Mozilla Sponsorship:** Gained Mozilla''s support in 2009, leading to increased development.',\n    '*   **Popularity:** Ranked as a top programming language since 2016.',\n    '*   **Community:** A large and active community contributes to the language''s development and ecosystem.',\n    '*   **Evolving:** Rust continually evolves with new features and improvements.',\n    '',\n    'Topic: Rust''s Memory Management Approach',\n    '',\n    'Rust''s approach to memory management differs significantly from that of garbage-collected languages and languages like C and C++. It avoids the overhead of garbage collection while providing the safety of managed memory. Rust achieves this through the concepts of ownership, borrowing, and lifetimes. Ownership dictates which part of the code is responsible for a given piece of memory. Borrowing allows multiple parts of the code to access the same memory, but under strict rules that prevent data races and ensure memory safety. Lifetimes help to track the duration of references, ensuring that they are valid for as long as they are being used.',\n    '',\n    '*   **Ownership:** Each value has a single owner, ensuring exclusive access.',\n    '*   **Borrowing:** Allows multiple parts of the code to access the same memory, but under strict rules to prevent data races.',\n    '*   **Lifetimes:** Track the validity of references to ensure they are not used after the memory they point to has been freed.',\n    '*   **Compile-Time Safety:** Memory management checks are performed at compile time.',\n    '*   **No Garbage Collector:** Ensures low overhead and predictability.',\n    '',\n    'Topic: Immutability and Memory Allocation in Rust',\n    '',\n    'By default, variables in Rust are immutable. This means that once a value is assigned to a variable, it cannot be changed. This immutability is a cornerstone of Rust''s memory safety, as it eliminates data races that can occur when multiple parts of a program try to modify the same memory location simultaneously. Immutable values, and objects with a known size at compile time, are typically stored on the stack, which is a fast and efficient memory region. Mutable values, as well as objects whose size is not known until runtime, are stored on the heap, a region of memory that is managed dynamically.',\n    '',\n    '*   **Immutability by Default:** Enhances safety and predictability.',\n    '*   **Stack Memory:** Used for immutable values, efficient due to minimal overhead.',\n    '*   **Heap Memory:** Used for mutable values and objects with unknown sizes at compile time.',\n    '*   **Compile-Time Allocation:** Memory allocation is often determined at compile time, optimizing performance.',\n    '*   **Automatic Dropping:** Memory is automatically freed when a variable goes out of scope, preventing memory leaks.'\n)\n\n# Replace the text of the first 49 paragraphs (1-based indices 1..49) in\n# place, preserving each paragraph's own paragraph mark. Setting Range.Text\n# on a paragraph's Range replaces only the run content before the mark; an\n# already-blank target paragraph is left untouched so it stays a bare <w:p/>\n# instead of gaining an empty run.\nfor ($i = 0; $i -lt $newTexts.Count; $i++) {\n    $p = $d.Paragraphs.Item($i + 1)\n    $newText = $newTexts[$i]\n    $currentText = $p.Range.Text.TrimEnd([char]13)\n    if ($newText -ne '' -or $currentText -ne '') {\n        $p.Range.Text = $newText\n    }\n}\n\n# The old document had a trailing \"Topic: Using NPM Modules\" section (10\n# paragraphs: #50 blank, #51 topic, #52 blank, #53 body, #54 blank, #55-59\n# bullets, #60 trailing blank) that the new document no longer has. Paragraph\n# #50 (the blank separator before that section) should remain as the\n# document's final (empty) paragraph, so repeatedly delete paragraph #50 --\n# each deletion shifts the following paragraphs up by one, and Word will not\n# delete the very last paragraph of the body, so this converges on removing\n# exactly the 10 obsolete paragraphs.\nfor ($i = 0; $i -lt 10; $i++) {\n    $p = $d.Paragraphs.Item(50)\n    $p.Range.Delete()\n}\n"}
